# This document contains a date line followed by a 20x5 table of
# addition/subtraction equations. Every equation text in the document is
# unique except for one pair of cells that happen to share identical text
# ("28+45=73"); that duplicate is handled by replacing one occurrence at a
# time (wdReplaceOne) in document order rather than a single ReplaceAll, so
# each of the two cells ends up with its own distinct new value.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-25 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-26 Thursday", 2) | Out-Null
$d.Content.Find.Execute("10+1=11", $true, $false, $false, $false, $false, $true, 1, $false, "51+30=81", 2) | Out-Null
$d.Content.Find.Execute("22+29=51", $true, $false, $false, $false, $false, $true, 1, $false, "27-23=4", 2) | Out-Null
$d.Content.Find.Execute("69-38=31", $true, $false, $false, $false, $false, $true, 1, $false, "36+11=47", 2) | Out-Null
$d.Content.Find.Execute("4+78=82", $true, $false, $false, $false, $false, $true, 1, $false, "24-16=8", 2) | Out-Null
$d.Content.Find.Execute("96-84=12", $true, $false, $false, $false, $false, $true, 1, $false, "11+18=29", 2) | Out-Null
$d.Content.Find.Execute("54-32=22", $true, $false, $false, $false, $false, $true, 1, $false, "15+19=34", 2) | Out-Null
$d.Content.Find.Execute("50-19=31", $true, $false, $false, $false, $false, $true, 1, $false, "0+51=51", 2) | Out-Null
$d.Content.Find.Execute("87-74=13", $true, $false, $false, $false, $false, $true, 1, $false, "75-51=24", 2) | Out-Null
$d.Content.Find.Execute("72+10=82", $true, $false, $false, $false, $false, $true, 1, $false, "24+9=33", 2) | Out-Null
$d.Content.Find.Execute("37-17=20", $true, $false, $false, $false, $false, $true, 1, $false, "27+60=87", 2) | Out-Null
$d.Content.Find.Execute("65+24=89", $true, $false, $false, $false, $false, $true, 1, $false, "96-22=74", 2) | Out-Null
$d.Content.Find.Execute("67-41=26", $true, $false, $false, $false, $false, $true, 1, $false, "75+9=84", 2) | Out-Null
$d.Content.Find.Execute("35-7=28", $true, $false, $false, $false, $false, $true, 1, $false, "0+61=61", 2) | Out-Null
$d.Content.Find.Execute("98-6=92", $true, $false, $false, $false, $false, $true, 1, $false, "59+13=72", 2) | Out-Null
$d.Content.Find.Execute("45-6=39", $true, $false, $false, $false, $false, $true, 1, $false, "9+56=65", 2) | Out-Null
$d.Content.Find.Execute("97-73=24", $true, $false, $false, $false, $false, $true, 1, $false, "12+82=94", 2) | Out-Null
$d.Content.Find.Execute("33+62=95", $true, $false, $false, $false, $false, $true, 1, $false, "41-8=33", 2) | Out-Null
$d.Content.Find.Execute("48-25=23", $true, $false, $false, $false, $false, $true, 1, $false, "43-26=17", 2) | Out-Null
$d.Content.Find.Execute("43-19=24", $true, $false, $false, $false, $false, $true, 1, $false, "11+13=24", 2) | Out-Null
$d.Content.Find.Execute("19+78=97", $true, $false, $false, $false, $false, $true, 1, $false, "5+67=72", 2) | Out-Null
$d.Content.Find.Execute("41-7=34", $true, $false, $false, $false, $false, $true, 1, $false, "45-44=1", 2) | Out-Null
$d.Content.Find.Execute("93-16=77", $true, $false, $false, $false, $false, $true, 1, $false, "67-53=14", 2) | Out-Null
$d.Content.Find.Execute("55+19=74", $true, $false, $false, $false, $false, $true, 1, $false, "21+58=79", 2) | Out-Null
$d.Content.Find.Execute("28+45=73", $true, $false, $false, $false, $false, $true, 1, $false, "43-36=7", 1) | Out-Null
$d.Content.Find.Execute("23+44=67", $true, $false, $false, $false, $false, $true, 1, $false, "86-4=82", 2) | Out-Null
$d.Content.Find.Execute("86-58=28", $true, $false, $false, $false, $false, $true, 1, $false, "86-25=61", 2) | Out-Null
$d.Content.Find.Execute("44-21=23", $true, $false, $false, $false, $false, $true, 1, $false, "80+7=87", 2) | Out-Null
$d.Content.Find.Execute("42-36=6", $true, $false, $false, $false, $false, $true, 1, $false, "74-33=41", 2) | Out-Null
$d.Content.Find.Execute("10-0=10", $true, $false, $false, $false, $false, $true, 1, $false, "2+78=80", 2) | Out-Null
$d.Content.Find.Execute("39+15=54", $true, $false, $false, $false, $false, $true, 1, $false, "44-44=0", 2) | Out-Null
$d.Content.Find.Execute("95-71=24", $true, $false, $false, $false, $false, $true, 1, $false, "10+54=64", 2) | Out-Null
$d.Content.Find.Execute("87-66=21", $true, $false, $false, $false, $false, $true, 1, $false, "53+39=92", 2) | Out-Null
$d.Content.Find.Execute("53-30=23", $true, $false, $false, $false, $false, $true, 1, $false, "91-10=81", 2) | Out-Null
$d.Content.Find.Execute("35-25=10", $true, $false, $false, $false, $false, $true, 1, $false, "73+11=84", 2) | Out-Null
$d.Content.Find.Execute("57+33=90", $true, $false, $false, $false, $false, $true, 1, $false, "56-25=31", 2) | Out-Null
$d.Content.Find.Execute("67-10=57", $true, $false, $false, $false, $false, $true, 1, $false, "49+48=97", 2) | Out-Null
$d.Content.Find.Execute("44+1=45", $true, $false, $false, $false, $false, $true, 1, $false, "15+7=22", 2) | Out-Null
$d.Content.Find.Execute("27+41=68", $true, $false, $false, $false, $false, $true, 1, $false, "27+68=95", 2) | Out-Null
$d.Content.Find.Execute("89-78=11", $true, $false, $false, $false, $false, $true, 1, $false, "36+16=52", 2) | Out-Null
$d.Content.Find.Execute("21+67=88", $true, $false, $false, $false, $false, $true, 1, $false, "17-7=10", 2) | Out-Null
$d.Content.Find.Execute("23+19=42", $true, $false, $false, $false, $false, $true, 1, $false, "14+44=58", 2) | Out-Null
$d.Content.Find.Execute("75-26=49", $true, $false, $false, $false, $false, $true, 1, $false, "44-18=26", 2) | Out-Null
$d.Content.Find.Execute("10+13=23", $true, $false, $false, $false, $false, $true, 1, $false, "80-12=68", 2) | Out-Null
$d.Content.Find.Execute("15+17=32", $true, $false, $false, $false, $false, $true, 1, $false, "24+47=71", 2) | Out-Null
$d.Content.Find.Execute("34-29=5", $true, $false, $false, $false, $false, $true, 1, $false, "12+31=43", 2) | Out-Null
$d.Content.Find.Execute("78-66=12", $true, $false, $false, $false, $false, $true, 1, $false, "84-59=25", 2) | Out-Null
$d.Content.Find.Execute("36+15=51", $true, $false, $false, $false, $false, $true, 1, $false, "7+74=81", 2) | Out-Null
$d.Content.Find.Execute("12+27=39", $true, $false, $false, $false, $false, $true, 1, $false, "73+8=81", 2) | Out-Null
$d.Content.Find.Execute("76-9=67", $true, $false, $false, $false, $false, $true, 1, $false, "50+33=83", 2) | Out-Null
$d.Content.Find.Execute("37-4=33", $true, $false, $false, $false, $false, $true, 1, $false, "72-61=11", 2) | Out-Null
$d.Content.Find.Execute("44-30=14", $true, $false, $false, $false, $false, $true, 1, $false, "54+12=66", 2) | Out-Null
$d.Content.Find.Execute("50+32=82", $true, $false, $false, $false, $false, $true, 1, $false, "92-84=8", 2) | Out-Null
$d.Content.Find.Execute("83-60=23", $true, $false, $false, $false, $false, $true, 1, $false, "91-37=54", 2) | Out-Null
$d.Content.Find.Execute("87-63=24", $true, $false, $false, $false, $false, $true, 1, $false, "72+7=79", 2) | Out-Null
$d.Content.Find.Execute("39-18=21", $true, $false, $false, $false, $false, $true, 1, $false, "57-52=5", 2) | Out-Null
$d.Content.Find.Execute("52-12=40", $true, $false, $false, $false, $false, $true, 1, $false, "64-60=4", 2) | Out-Null
$d.Content.Find.Execute("81-77=4", $true, $false, $false, $false, $false, $true, 1, $false, "35+5=40", 2) | Out-Null
$d.Content.Find.Execute("23-15=8", $true, $false, $false, $false, $false, $true, 1, $false, "87-42=45", 2) | Out-Null
$d.Content.Find.Execute("67-35=32", $true, $false, $false, $false, $false, $true, 1, $false, "89-27=62", 2) | Out-Null
$d.Content.Find.Execute("38-19=19", $true, $false, $false, $false, $false, $true, 1, $false, "74-57=17", 2) | Out-Null
$d.Content.Find.Execute("70-55=15", $true, $false, $false, $false, $false, $true, 1, $false, "31+2=33", 2) | Out-Null
$d.Content.Find.Execute("37-32=5", $true, $false, $false, $false, $false, $true, 1, $false, "97-49=48", 2) | Out-Null
$d.Content.Find.Execute("29+58=87", $true, $false, $false, $false, $false, $true, 1, $false, "95-53=42", 2) | Out-Null
$d.Content.Find.Execute("94-24=70", $true, $false, $false, $false, $false, $true, 1, $false, "27-10=17", 2) | Out-Null
$d.Content.Find.Execute("65+16=81", $true, $false, $false, $false, $false, $true, 1, $false, "98-45=53", 2) | Out-Null
$d.Content.Find.Execute("85-16=69", $true, $false, $false, $false, $false, $true, 1, $false, "60-33=27", 2) | Out-Null
$d.Content.Find.Execute("64+26=90", $true, $false, $false, $false, $false, $true, 1, $false, "14+72=86", 2) | Out-Null
$d.Content.Find.Execute("75-38=37", $true, $false, $false, $false, $false, $true, 1, $false, "59-37=22", 2) | Out-Null
$d.Content.Find.Execute("37+49=86", $true, $false, $false, $false, $false, $true, 1, $false, "79-64=15", 2) | Out-Null
$d.Content.Find.Execute("75+15=90", $true, $false, $false, $false, $false, $true, 1, $false, "66-40=26", 2) | Out-Null
$d.Content.Find.Execute("53-32=21", $true, $false, $false, $false, $false, $true, 1, $false, "76-36=40", 2) | Out-Null
$d.Content.Find.Execute("15+76=91", $true, $false, $false, $false, $false, $true, 1, $false, "49+27=76", 2) | Out-Null
$d.Content.Find.Execute("47+15=62", $true, $false, $false, $false, $false, $true, 1, $false, "1+63=64", 2) | Out-Null
$d.Content.Find.Execute("36-25=11", $true, $false, $false, $false, $false, $true, 1, $false, "16+39=55", 2) | Out-Null
$d.Content.Find.Execute("82-35=47", $true, $false, $false, $false, $false, $true, 1, $false, "93-34=59", 2) | Out-Null
$d.Content.Find.Execute("37+38=75", $true, $false, $false, $false, $false, $true, 1, $false, "61-59=2", 2) | Out-Null
$d.Content.Find.Execute("28+45=73", $true, $false, $false, $false, $false, $true, 1, $false, "57+4=61", 1) | Out-Null
$d.Content.Find.Execute("10+23=33", $true, $false, $false, $false, $false, $true, 1, $false, "35+54=89", 2) | Out-Null
$d.Content.Find.Execute("73-6=67", $true, $false, $false, $false, $false, $true, 1, $false, "32-12=20", 2) | Out-Null
$d.Content.Find.Execute("97-23=74", $true, $false, $false, $false, $false, $true, 1, $false, "77-4=73", 2) | Out-Null
$d.Content.Find.Execute("91-4=87", $true, $false, $false, $false, $false, $true, 1, $false, "9+57=66", 2) | Out-Null
$d.Content.Find.Execute("64-4=60", $true, $false, $false, $false, $false, $true, 1, $false, "1+16=17", 2) | Out-Null
$d.Content.Find.Execute("58-3=55", $true, $false, $false, $false, $false, $true, 1, $false, "60-49=11", 2) | Out-Null
$d.Content.Find.Execute("88-23=65", $true, $false, $false, $false, $false, $true, 1, $false, "78-32=46", 2) | Out-Null
$d.Content.Find.Execute("24+35=59", $true, $false, $false, $false, $false, $true, 1, $false, "11+65=76", 2) | Out-Null
$d.Content.Find.Execute("95-44=51", $true, $false, $false, $false, $false, $true, 1, $false, "15+2=17", 2) | Out-Null
$d.Content.Find.Execute("47+51=98", $true, $false, $false, $false, $false, $true, 1, $false, "2+52=54", 2) | Out-Null
$d.Content.Find.Execute("46-38=8", $true, $false, $false, $false, $false, $true, 1, $false, "4+29=33", 2) | Out-Null
$d.Content.Find.Execute("36+55=91", $true, $false, $false, $false, $false, $true, 1, $false, "30-6=24", 2) | Out-Null
$d.Content.Find.Execute("62-11=51", $true, $false, $false, $false, $false, $true, 1, $false, "9+35=44", 2) | Out-Null
$d.Content.Find.Execute("4+13=17", $true, $false, $false, $false, $false, $true, 1, $false, "6+0=6", 2) | Out-Null
$d.Content.Find.Execute("10+35=45", $true, $false, $false, $false, $false, $true, 1, $false, "55-22=33", 2) | Out-Null
$d.Content.Find.Execute("33+8=41", $true, $false, $false, $false, $false, $true, 1, $false, "48-45=3", 2) | Out-Null
$d.Content.Find.Execute("87-12=75", $true, $false, $false, $false, $false, $true, 1, $false, "93-46=47", 2) | Out-Null
$d.Content.Find.Execute("34-22=12", $true, $false, $false, $false, $false, $true, 1, $false, "43-6=37", 2) | Out-Null
$d.Content.Find.Execute("97-71=26", $true, $false, $false, $false, $false, $true, 1, $false, "89-60=29", 2) | Out-Null
$d.Content.Find.Execute("58+16=74", $true, $false, $false, $false, $false, $true, 1, $false, "3+15=18", 2) | Out-Null
$d.Content.Find.Execute("93+0=93", $true, $false, $false, $false, $false, $true, 1, $false, "42+28=70", 2) | Out-Null
$d.Content.Find.Execute("94-73=21", $true, $false, $false, $false, $false, $true, 1, $false, "41+58=99", 2) | Out-Null
$d.Content.Find.Execute("39-25=14", $true, $false, $false, $false, $false, $true, 1, $false, "38-30=8", 2) | Out-Null
